$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Updated account-statement detail rows (B15:J22 table).
# Employee identity (col C doc number, col D name), period (col E),
# "Salario Basico" (col F) and "Valor Mora" (col G) have all been
# refreshed with the new database values.

$ws.Range("C16").Value = "1044938732"
$ws.Range("D16").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E16").Value = "2211"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1707149

$ws.Range("C17").Value = "1128059537"
$ws.Range("D17").Value = "LIZZETH PAOLA OSPINO GONZALEZ"
$ws.Range("E17").Value = "2212"
$ws.Range("F17").Value = 88000
$ws.Range("G17").Value = 3029734

$ws.Range("C18").Value = "1044938732"
$ws.Range("D18").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E18").Value = "2212"
$ws.Range("F18").Value = 48000
$ws.Range("G18").Value = 1707149

$ws.Range("C19").Value = "1044938732"
$ws.Range("D19").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 48000
$ws.Range("G19").Value = 1707149

$ws.Range("C20").Value = "1128059537"
$ws.Range("D20").Value = "LIZZETH PAOLA OSPINO GONZALEZ"
$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 88000
$ws.Range("G20").Value = 3029734

$ws.Range("C21").Value = "1044938732"
$ws.Range("D21").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E21").Value = "2302"
$ws.Range("F21").Value = 48000
$ws.Range("G21").Value = 1707149

$ws.Range("C22").Value = "1044938732"
$ws.Range("D22").Value = "CATERIN ALEZANDRA MAZA PALOMINO"
$ws.Range("E22").Value = "2304"
$ws.Range("F22").Value = 48000
$ws.Range("G22").Value = 1707149

$wb.Save()
